# Applies the "Improving the custom plots" edit to the velocity_tuning
# plot-definition sheet:
#   * Rename the plot title in B3 from "Velocity Tuning" to "VelocityTuning"
#   * Insert a new row into the "PN" plot group (Plot 2) wiring up the
#     PSC/TPX -> PN_{DES} desired-north-position channel
#   * Insert a new row into the "PE" plot group (Plot 4) wiring up the
#     PSC/TPY -> PE_{DES} desired-east-position channel

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the plot title -------------------------------------------------
$ws.Range("B3").Value = "VelocityTuning"

# --- Insert the two new data rows ------------------------------------------
# Row 11 becomes the 3rd row of the "PN" plot group (pushes the old rows
# 11-15 down to 12-16). Row 17 (inserted after the data that is now in row
# 16) becomes the 3rd row of the "PE" plot group.
$ws.Rows("11:11").Insert()
$ws.Rows("17:17").Insert()

# --- Fill in the new "PN_{DES}" row (row 11) --------------------------------
$ws.Range("A11").Value2 = 2
$ws.Range("B11").Value2 = 2
$ws.Range("C11").Value2 = 1
$ws.Range("D11").Value = "S"
$ws.Range("E11").Value = "Time [ s ]"
$ws.Range("F11").Value = "PN~[~m~]"
$ws.Range("G11").Value = "Vertical"
$ws.Range("K11").Value = "PSC/TPX"
$ws.Range("O11").Value2 = 1
$ws.Range("R11").Value = "PN_{DES}"
$ws.Range("S11").Value = "m"

# --- Fill in the new "PE_{DES}" row (row 17) --------------------------------
$ws.Range("A17").Value2 = 4
$ws.Range("B17").Value2 = 2
$ws.Range("C17").Value2 = 2
$ws.Range("D17").Value = "S"
$ws.Range("E17").Value = "Time [ s ]"
$ws.Range("F17").Value = "PE~[~m~]"
$ws.Range("G17").Value = "Vertical"
$ws.Range("K17").Value = "PSC/TPY"
$ws.Range("O17").Value2 = 1
$ws.Range("R17").Value = "PE_{DES}"
$ws.Range("S17").Value = "m"

# --- Tidy up the sheet view/selection --------------------------------------
$ws.Range("B3").Select()
